# ----------------------------------------------------------------------
# Applies two small text corrections to the resume:
#   1) "ISO 27002" -> "ISO 27001" in the "Misc:" skills line (the digit
#      "2" at the end of "27002" becomes "1"; Word splits the run around
#      the single changed character).
#   2) "ArcSight EMS deployment" -> "ArcSight ESM deployment" in the
#      bullet about the security-architect role.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) ISO 27002 -> ISO 27001 -------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("ISO 27002", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $find.Start
    $matchEnd   = $find.End

    # Only the last character ("2" -> "1") actually changes; isolate it as
    # its own Range so the edit lands on a single character, the same way
    # Word would if you selected that character and retyped it.
    $digitRange = $d.Range($matchEnd - 1, $matchEnd)
    $digitRange.Text = "1"

    # Nudge the newly-typed character's formatting so Word keeps it as its
    # own run instead of silently re-merging it with its neighbours (the
    # same run-splitting you see after a manual retype in the real app).
    $newCharRange = $d.Range($matchEnd - 1, $matchEnd)
    $newCharRange.Font.Bold = 9999999
}

# --- 2) ArcSight EMS -> ArcSight ESM -------------------------------------------
$d.Content.Find.Execute("ArcSight EMS deployment", $true, $false, $false, $false, $false, $true, 1, $false, "ArcSight ESM deployment", 2) | Out-Null
